$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

$data = @{
  2 = @(1.02,1.039853672188411,1.040342139962946,1.054863011287509,1.061339771427892,1.033350078465052,1.044943438887744,1.043124760710345,1.057605051498156,1.064064099261173,1.046427378427653)
  3 = @(1.02,1.041587208595509,1.041643908738153,1.056529615243884,1.063178780174586,1.033713665860834,1.046318436330593,1.044235640677261,1.059082805258739,1.065715130068972,1.047804328524584)
  4 = @(1.02,1.042705623798089,1.042483173613096,1.057606042126113,1.064366712817129,1.033946269798782,1.047204592128104,1.044950852837978,1.06003656787782,1.066780996014021,1.048691742764971)
  5 = @(1.02,1.043175031341449,1.042835277002884,1.058058112018526,1.064865648965061,1.034043423168588,1.047576292156386,1.045250675747798,1.060436957275962,1.067228510619541,1.049063970649928)
  6 = @(1.02,1.043253801981885,1.042894354586581,1.058133989957227,1.064949395369268,1.034059698618619,1.0476386533309,1.045300967641504,1.060504151210422,1.067303616946397,1.049126420384463)
  7 = @(1.02,1.042711899065532,1.04248788126438,1.057612084495133,1.064373381455477,1.033947570450114,1.047209562085345,1.044954862421633,1.060041920132697,1.066786977968226,1.048696719780122)
  8 = @(1.02,1.040440221352581,1.040782720168159,1.055426665417039,1.06196170223064,1.033473507906379,1.045408870154294,1.043500939499801,1.058104978623481,1.064622590855328,1.046893470659986)
  9 = @(1.02,1.03641127011902,1.037754044423573,1.051559929761728,1.057695799338938,1.032617592876892,1.042208005196337,1.040910910879759,1.05467257074915,1.060789189848955,1.043688060106645)
  10 = @(1.02,1.033706861105658,1.035718142243306,1.048970669416594,1.054840011531996,1.032032917750234,1.040054613650803,1.039164748164592,1.052370556543384,1.058219620405827,1.041531610498119)
  11 = @(1.02,1.032531227601059,1.034832441777376,1.047846588621185,1.053600400097783,1.031776355849216,1.039117375636874,1.038403879068888,1.051370323418807,1.057103459340675,1.040593041498707)
  12 = @(1.02,1.032093833775211,1.034502817838693,1.047428600727244,1.053139478623869,1.031680542643316,1.038768506572264,1.038120529566974,1.050998259352207,1.056688321527985,1.04024367700001)
  13 = @(1.02,1.032187688728116,1.034573552256331,1.047518281370572,1.05323836961865,1.031701118270294,1.038843373773832,1.038181342178333,1.051078092721388,1.056777394975582,1.040318650521592)
  14 = @(1.02,1.032495087110581,1.034805207997412,1.047812046973577,1.053562309973721,1.031768446427443,1.03908855313711,1.038380472250754,1.051339579447889,1.057069155154592,1.040564178067688)
  15 = @(1.02,1.032684390583116,1.034947854158715,1.047992984971109,1.053761836875565,1.031809861179964,1.039239518099199,1.038503065910099,1.051500618969858,1.057248845382465,1.040715357417339)
  16 = @(1.02,1.033784785348424,1.035776834893074,1.049045208163379,1.054922214919825,1.03204987300659,1.040116712516619,1.0392151430114,1.052436864824556,1.058293620721593,1.041593797551447)
  17 = @(1.02,1.034473787185125,1.036295714605734,1.049704448266934,1.055649263538654,1.032199514104512,1.040665656222361,1.039660524676109,1.053023213862494,1.058948026218487,1.042143520820216)
  18 = @(1.02,1.034875227136422,1.036597969676807,1.050088691829338,1.056073046975498,1.03228647011054,1.040985383196474,1.039919848956679,1.053364889907556,1.0593293913486,1.042463701843305)
  19 = @(1.02,1.035012033192366,1.036700963616239,1.05021966195001,1.056217497248311,1.032316064522076,1.041094323823789,1.040008194298782,1.053481336852206,1.059459370093583,1.042572797178807)
  20 = @(1.019999999999999,1.034399909785023,1.03624008505613,1.049633747126991,1.055571288409509,1.032183492900718,1.040606807716645,1.039612787036046,1.052960338551211,1.058877849839588,1.042084588742872)
  21 = @(1.02,1.03240458573582,1.034737008800754,1.047725552997562,1.053466930863064,1.031748634210373,1.039016374397104,1.038321853627065,1.051262592959724,1.056983254260503,1.040491896825592)
  22 = @(1.02,1.03114592190842,1.033788282527395,1.046523159256065,1.052141082632783,1.031472241003616,1.038012133339892,1.037505969476004,1.05019206176029,1.055788879886606,1.039486229630891)
  23 = @(1.02,1.031813560805182,1.034291573775331,1.047160826363289,1.052844207216228,1.03161904639655,1.038544910632227,1.037938889559683,1.050759868750357,1.056422346247545,1.040019763528091)
  24 = @(1.02,1.034433293171144,1.036265222886847,1.049665694793922,1.055606522915115,1.03219073320381,1.040633400252473,1.03963435903984,1.05298875022183,1.058909560570288,1.04211121904315)
  25 = @(1.02,1.037456028192323,1.038539938337251,1.052561526634529,1.058800656329875,1.03284132839671,1.043038881538278,1.041583881597189,1.05556229445121,1.061782613241197,1.044520116388294)
}

foreach ($r in $data.Keys) {
  $vals = $data[$r]
  for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + $r).Value = $vals[$i]
  }
}

Write-Output "done"